$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen columns D and E (closest value Excel's pixel-rounded
# ColumnWidth can store to the target 33.7109375 / 27.7109375
# character widths; COM ColumnWidth always snaps to a whole
# pixel for the sheet's font, so an exact hit isn't addressable).
$ws.Columns.Item(4).ColumnWidth = 32.8
$ws.Columns.Item(5).ColumnWidth = 26.8

# Frame the two newly-added rows (16 and 17) with a thin box
# border on every cell A:E, matching the new border/cellXfs
# records added to styles.xml.
$ws.Range("A16:E17").Borders.LineStyle = 1

# Move the active selection to G17 (previously F12).
$ws.Range("G17").Select()
